$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")

# Give new row 16 the same formatting as the existing row 15 (plain data row)
# before any values move, so the style table isn't enlarged with an extra xf.
$wsMeta.Range("A15:B15").Copy()
$wsMeta.Range("A16:B16").PasteSpecial(-4122)

# Update Version
$wsMeta.Range("B3").Value = "0.1.7"

# Update Status
$wsMeta.Range("B6").Value = "draft"

# Update Date
$wsMeta.Range("B8").Value = "2024-11-22T12:33:30-06:00"

# Update first Contact (organization) text
$wsMeta.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Update second Contact (person) text - row already existed as a duplicate placeholder
$wsMeta.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Shift Description/Purpose/Copyright/Immutable rows down by one (to rows 13-16) and
# insert the new Jurisdiction row at 12 - values are written directly (bottom-up so
# no row's content is clobbered before it has been moved along).
$wsMeta.Range("A16").Value = "Immutable"
$wsMeta.Range("B16").Value = "BooleanType[null]"

$wsMeta.Range("A15").Value = "Copyright"
$wsMeta.Range("B15").Value = ""

$wsMeta.Range("A14").Value = "Purpose"
$wsMeta.Range("B14").Value = ""

$wsMeta.Range("A13").Value = "Description"
$wsMeta.Range("B13").Value = "RxNorm codes for Treosulfan"

$wsMeta.Range("A12").Value = "Jurisdiction"
$wsMeta.Range("B12").Value = ""

$wb.Save()
